# Update "想去人数" (number of interested attendees) figures on the
# "展览" (Exhibition) and "全部类型" (All types) worksheets to reflect
# the newly generated data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 15225
$ws1.Range("F4").Value = 711
$ws1.Range("F5").Value = 250
$ws1.Range("F6").Value = 642
$ws1.Range("F7").Value = 1619

# --- Sheet "全部类型" ------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 15225
$ws4.Range("F4").Value = 711
$ws4.Range("F5").Value = 250
$ws4.Range("F8").Value = 642
$ws4.Range("F9").Value = 1619
